$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update countries COVID figures (refresh of the source data) and bump the
# "datos actualizados" timestamp in the title cell.

# Update "datos actualizados" timestamp title
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 13:47"

# Row 4: Estados Unidos (updated totals)
$ws.Range("B4").Value = 7601182
$ws.Range("C4").Value = 336
$ws.Range("D4").Value = 4818768
$ws.Range("E4").Value = 2568134
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 214280

# Row 16: Chile -> Iran
$ws.Range("A16").Value = "Iran"
$ws.Range("B16").Value = 471772
$ws.Range("C16").Value = 3653
$ws.Range("D16").Value = 389966
$ws.Range("E16").Value = 54849
$ws.Range("G16").Value = 211
$ws.Range("H16").Value = 26957

# Row 17: Iran -> Chile
$ws.Range("A17").Value = "Chile"
$ws.Range("B17").Value = 468471
$ws.Range("D17").Value = 440881
$ws.Range("E17").Value = 14671
$ws.Range("H17").Value = 12919

# Row 47: China -> Nepal
$ws.Range("A47").Value = "Nepal"
$ws.Range("B47").Value = 86823
$ws.Range("C47").Value = 2253
$ws.Range("D47").Value = 64069
$ws.Range("E47").Value = 22219
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 535

# Row 48: Japon -> China
$ws.Range("A48").Value = "China"
$ws.Range("B48").Value = 85450
$ws.Range("C48").Value = 16
$ws.Range("D48").Value = 80621
$ws.Range("E48").Value = 195
$ws.Range("H48").Value = 4634

# Row 49: Nepal -> Japon
$ws.Range("A49").Value = "Japon"
$ws.Range("B49").Value = 84768
$ws.Range("D49").Value = 77807
$ws.Range("E49").Value = 5371
$ws.Range("H49").Value = 1590

# Row 67: Ghana (updated totals)
$ws.Range("B67").Value = 46829
$ws.Range("C67").Value = 26
$ws.Range("D67").Value = 46060
$ws.Range("E67").Value = 466

# Row 83: Bulgaria -> Tunez
$ws.Range("A83").Value = "Tunez"
$ws.Range("B83").Value = 22230
$ws.Range("C83").Value = 1286
$ws.Range("D83").Value = 5032
$ws.Range("E83").Value = 16877
$ws.Range("G83").Value = 45
$ws.Range("H83").Value = 321

# Row 84: Tunez -> Bulgaria
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 21518
$ws.Range("D84").Value = 14984
$ws.Range("E84").Value = 5693
$ws.Range("H84").Value = 841

# Row 90: Madagascar (updated totals)
$ws.Range("B90").Value = 16558
$ws.Range("C90").Value = 29
$ws.Range("D90").Value = 15486
$ws.Range("E90").Value = 840

# Row 144: Malta (updated totals)
$ws.Range("B144").Value = 3270
$ws.Range("C144").Value = 66
$ws.Range("D144").Value = 2758
$ws.Range("E144").Value = 473
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 39

# Row 163: Lesoto (updated totals)
$ws.Range("B163").Value = 1683
$ws.Range("C163").Value = 3
$ws.Range("E163").Value = 718
